$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.620.01'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '2.473.25'

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = '''317.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.42%  '

$ws.Range("D6").Value = '''92.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.54%  '

$ws.Range("E7").Value = '  +1.45%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").Value = '''0.515'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.74%  '

$ws.Range("B10").Value = 'Avalanche'
$ws.Range("C10").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D10").Value = '''33.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.21%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '''0.0860'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.68%  '

$ws.Range("E12").Value = '  +0.81%  '

$ws.Range("D13").Value = '2.854.77'

$ws.Range("D14").Value = '''6.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.64%  '

$ws.Range("D15").Value = '''15.64'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.96%  '

$ws.Range("D16").Value = '2.477.05'
$ws.Range("E16").Value = '  -1.06%  '

$ws.Range("D17").Value = '''0.790'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.04%  '

$ws.Range("D18").Value = '41.585.52'
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("D19").Value = '''6.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.56%  '

$ws.Range("D20").Value = '0.0₃0952'
$ws.Range("E20").Value = '  +1.40%  '

$ws.Range("D21").Value = '''71.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.56%  '

$ws.Range("D22").Value = '''11.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.72%  '

$ws.Range("D23").Value = '''240.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.51%  '

$ws.Range("E24").Value = '  +1.94%  '

$ws.Range("D25").Value = '''1.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.56%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").Value = '''24.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.54%  '

$ws.Range("E28").Value = '  +3.91%  '

$ws.Range("D29").Value = '''9.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.46%  '

$ws.Range("D30").Value = '''36.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.02%  '

$ws.Range("D31").Value = '''160.21'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.28%  '

$ws.Range("D32").Value = '''5.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.41%  '

$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("D34").Value = '''0.0771'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.99%  '

$ws.Range("E35").Value = '  +0.15%  '

$ws.Range("D36").Value = '''17.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.23%  '

$ws.Range("E37").Value = '  -0.11%  '

$ws.Range("E38").Value = '  +0.77%  '

$ws.Range("E39").Value = '  +1.56%  '

$ws.Range("E40").Value = '  -2.94%  '

$ws.Range("E41").Value = '  -2.95%  '

$ws.Range("E42").Value = '  +2.03%  '

$ws.Range("D43").Value = '1.989.76'
$ws.Range("E43").Value = '  +0.50%  '

$ws.Range("D44").Value = '''19.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.28%  '

$ws.Range("E45").Value = '  +1.05%  '

$ws.Range("D46").Value = '''2.99'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.44%  '

$ws.Range("E47").Value = '  +2.76%  '

$ws.Range("D48").Value = '2.712.18'
$ws.Range("E48").Value = '  -0.28%  '

$ws.Range("D49").Value = '''97.64'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("D50").Value = '''67.30'
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").Value = '''73.57'
$ws.Range("D51").Style = "Normal"
